# The data-entry app this workbook feeds had a "Date of Exam" line that
# previously only captured a single date; it's been updated to handle
# multiple date matches on the same line, so the sample/dummy data needs
# a row exercising that case.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Patients")

$ws.Range("B3").Value = "``Date of Exam: 4-1-2023, April 1st, 2023, May 2nd, 2023, June 3rd, 2023, July 4th, 2023. Air pressure: 1234.5 mBar.'"

# Column B was sized to fit its longest entry; re-fit it now that B3 grew.
$ws.Columns.Item(2).ColumnWidth = 101

# Leave the selection where the editor left it after making this change.
$ws.Range("B15:B16").Select()
